# Prior inputs are now mean and std
# Rework Planilha1's solver-helper block (cols H:J) so that instead of
# storing raw "SD"/"obj" helper numbers and a sum-of-squares objective in
# column J, column H now derives the distribution's first shape parameter
# from the already-computed Mean (D) / Variance (E) / Std (F), and column I
# derives the second shape parameter from H - i.e. the prior's Param a/ Param b
# inputs become expressed as mean & std. The stray column J (objective) is
# removed entirely, and the B2/C2 inputs for the Beta row become 1.5 / 1.5.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Planilha1")
$ws2 = $wb.Worksheets.Item("Iterações")

# --- Planilha1, row 2 (Beta) inputs: Param a / Param b go from 3 / 2 to 1.5 / 1.5
$ws1.Range("B2").Value = 1.5
$ws1.Range("C2").Value = 1.5

# --- Header row: H1/I1 relabeled from "Mean"/"SD" to "a"/"b" (the derived
# shape parameters), matching the bold centered style used by the other
# headers; the stray J1 ("obj") header is cleared
$ws1.Range("H1").Value = "a"
$ws1.Range("H1").Font.Bold = $true
$ws1.Range("H1").HorizontalAlignment = -4108
$ws1.Range("I1").Value = "b"
$ws1.Range("I1").Font.Bold = $true
$ws1.Range("I1").HorizontalAlignment = -4108
$ws1.Range("J1").ClearContents()

# --- Row 2: H2/I2 become formulas derived from D2/E2; J2 (objective) cleared
$ws1.Range("H2").Formula = "=((D2^2)*(1-D2))/E2-D2"
$ws1.Range("I2").Formula = "=H2*(1-D2)/D2"
$ws1.Range("J2").ClearContents()

# --- Row 3: H3/I3 become formulas derived from D3/F3; J3 cleared
$ws1.Range("H3").Formula = "=(D3/F3)^2"
$ws1.Range("I3").Formula = "=D3/H3"
$ws1.Range("J3").ClearContents()

# --- Row 4: H4/I4 become formulas derived from D4/F4/H4; J4 cleared
$ws1.Range("H4").Formula = "=(D4/F4)^2+2"
$ws1.Range("I4").Formula = "=D4*(H4-1)"
$ws1.Range("J4").ClearContents()

# --- Row 5: H5/I5 become formulas derived from D5/F5/H5; J5 cleared
$ws1.Range("H5").Formula = "=D5-SQRT(3)*F5"
$ws1.Range("I5").Formula = "=2*D5-H5"
$ws1.Range("J5").ClearContents()

# --- Row 6: H6/I6 become formulas derived from D6/F6; J6 cleared
$ws1.Range("H6").Formula = "=D6"
$ws1.Range("I6").Formula = "=F6"
$ws1.Range("J6").ClearContents()

# --- The "theta"/"sigma_pi" labels used lower in the sheet (A11/A18) keep
# their text; no change needed there since the shared strings already read
# "theta" and "sigma_pi".

# --- Sheet view bookkeeping: Planilha1 becomes the active/selected sheet
# (was previously on "Iterações"); restore plain selections on both sheets.
[void]$ws2.Range("D2").Select()
[void]$ws1.Range("D2").Select()
[void]$ws1.Activate()
